$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 713929.0600000001
$ws.Range("I137").Value = 2271729.2
$ws.Range("K137").Value = 6815187.600000001
$ws.Range("M137").Value = -6812637.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6066
$ws.Range("I32").Value = 8118.12
$ws.Range("J32").Value = 3734.0454
$ws.Range("K32").Value = 8118.12
$ws.Range("L32").Value = 3734.0454
$ws.Range("M32").Value = -7831.12
$ws.Range("N32").Value = -4308.0454

$ws.Range("H61").Value = 1194.0646
$ws.Range("I61").Value = 1183.5652
$ws.Range("J61").Value = 1224.25
$ws.Range("K61").Value = 1183.5652
$ws.Range("L61").Value = 1224.25
$ws.Range("M61").Value = -971.5652
$ws.Range("N61").Value = -1648.25

$ws.Range("H74").Value = 314771.9
$ws.Range("I74").Value = 395376.34
$ws.Range("J74").Value = 3869.1428
$ws.Range("K74").Value = 395376.34
$ws.Range("L74").Value = 3869.1428
$ws.Range("M74").Value = -394502.34
$ws.Range("N74").Value = -5617.1428

$ws.Range("H77").Value = 314771.9
$ws.Range("I77").Value = 395376.34
$ws.Range("J77").Value = 3869.1428
$ws.Range("K77").Value = 1976881.7
$ws.Range("L77").Value = 19345.714
$ws.Range("M77").Value = -1972513.7
$ws.Range("N77").Value = -28081.714

$ws.Range("H88").Value = 16669111
$ws.Range("I88").Value = 33334082
$ws.Range("K88").Value = 33334082
$ws.Range("M88").Value = -33333676

$ws.Range("H91").Value = 16669111
$ws.Range("I91").Value = 33334082
$ws.Range("K91").Value = 33334082
$ws.Range("M91").Value = -33332678

$ws.Range("H122").Value = 3609.3713
$ws.Range("I122").Value = 3538.7083
$ws.Range("J122").Value = 3763.5454
$ws.Range("K122").Value = 10616.1249
$ws.Range("L122").Value = 11290.6362
$ws.Range("M122").Value = -8166.124899999999
$ws.Range("N122").Value = -16190.6362

$ws.Range("H136").Value = 1194.0646
$ws.Range("I136").Value = 1183.5652
$ws.Range("J136").Value = 1224.25
$ws.Range("K136").Value = 3550.6956
$ws.Range("L136").Value = 3672.75
$ws.Range("M136").Value = -1000.6956
$ws.Range("N136").Value = -8772.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2832.2856
$ws.Range("I86").Value = 2832.2856
$ws.Range("K86").Value = 2832.2856
$ws.Range("M86").Value = -1709.2856

$ws.Range("H89").Value = 2832.2856
$ws.Range("I89").Value = 2832.2856
$ws.Range("K89").Value = 14161.428
$ws.Range("M89").Value = -8545.428

$ws.Range("H134").Value = 3287.1482
$ws.Range("I134").Value = 1098.8
$ws.Range("J134").Value = 6022.5835
$ws.Range("K134").Value = 3296.4
$ws.Range("L134").Value = 18067.7505
$ws.Range("M134").Value = -761.3999999999996
$ws.Range("N134").Value = -23137.7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 200437.23
$ws.Range("I31").Value = 398233.72
$ws.Range("K31").Value = 398233.72
$ws.Range("M31").Value = -397938.72

$ws.Range("H34").Value = 200437.23
$ws.Range("I34").Value = 398233.72
$ws.Range("K34").Value = 398233.72
$ws.Range("M34").Value = -398031.72

$ws.Range("H58").Value = 2624.6858
$ws.Range("I58").Value = 1462.5
$ws.Range("K58").Value = 1462.5
$ws.Range("M58").Value = -1259.5

$ws.Range("H86").Value = 2754.389
$ws.Range("I86").Value = 2552.4614
$ws.Range("J86").Value = 3279.4
$ws.Range("K86").Value = 2552.4614
$ws.Range("L86").Value = 3279.4
$ws.Range("M86").Value = -1429.4614
$ws.Range("N86").Value = -5525.4

$ws.Range("H89").Value = 2754.389
$ws.Range("I89").Value = 2552.4614
$ws.Range("J89").Value = 3279.4
$ws.Range("K89").Value = 12762.307
$ws.Range("L89").Value = 16397
$ws.Range("M89").Value = -7146.307000000001
$ws.Range("N89").Value = -27629

$ws.Range("H132").Value = 2598.8462
$ws.Range("I132").Value = 2201.7273
$ws.Range("J132").Value = 4783
$ws.Range("K132").Value = 6605.1819
$ws.Range("L132").Value = 14349
$ws.Range("M132").Value = -4075.1819
$ws.Range("N132").Value = -19409

$ws.Range("H134").Value = 1277
$ws.Range("I134").Value = 856.89655
$ws.Range("K134").Value = 2570.68965
$ws.Range("M134").Value = -35.68965000000026

$ws.Range("H136").Value = 2624.6858
$ws.Range("I136").Value = 1462.5
$ws.Range("K136").Value = 4387.5
$ws.Range("M136").Value = -1837.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1567.6296
$ws.Range("I68").Value = 1218.2
$ws.Range("J68").Value = 1702.0256
$ws.Range("K68").Value = 3654.6
$ws.Range("L68").Value = 5106.0768
$ws.Range("M68").Value = -2843.6
$ws.Range("N68").Value = -6728.0768

$ws.Range("H71").Value = 1567.6296
$ws.Range("I71").Value = 1218.2
$ws.Range("J71").Value = 1702.0256
$ws.Range("K71").Value = 10963.8
$ws.Range("L71").Value = 15318.2304
$ws.Range("M71").Value = -6907.800000000001
$ws.Range("N71").Value = -23430.2304

$ws.Range("H107").Value = 6186310
$ws.Range("I107").Value = 472.76315
$ws.Range("J107").Value = 11652864
$ws.Range("K107").Value = 1418.28945
$ws.Range("L107").Value = 34958592
$ws.Range("M107").Value = 501.71055
$ws.Range("N107").Value = -34962432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6724.162
$ws.Range("I70").Value = 6017.3213
$ws.Range("J70").Value = 8923.223
$ws.Range("K70").Value = 6017.3213
$ws.Range("L70").Value = 8923.223
$ws.Range("M70").Value = -5747.3213
$ws.Range("N70").Value = -9463.223

$ws.Range("H73").Value = 6724.162
$ws.Range("I73").Value = 6017.3213
$ws.Range("J73").Value = 8923.223
$ws.Range("K73").Value = 6017.3213
$ws.Range("L73").Value = 8923.223
$ws.Range("M73").Value = -5081.3213
$ws.Range("N73").Value = -10795.223

$ws.Range("H126").Value = 3453.21
$ws.Range("I126").Value = 2807.8713
$ws.Range("J126").Value = 4959
$ws.Range("K126").Value = 8423.6139
$ws.Range("L126").Value = 14877
$ws.Range("M126").Value = -5953.6139
$ws.Range("N126").Value = -19817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2897.4333
$ws.Range("I7").Value = 1339.1052
$ws.Range("J7").Value = 5589.091
$ws.Range("K7").Value = 1339.1052
$ws.Range("L7").Value = 5589.091
$ws.Range("M7").Value = -1227.1052
$ws.Range("N7").Value = -5813.091

$ws.Range("H122").Value = 5407.8823
$ws.Range("I122").Value = 4595.6
$ws.Range("J122").Value = 11500
$ws.Range("K122").Value = 13786.8
$ws.Range("L122").Value = 34500
$ws.Range("M122").Value = -11336.8
$ws.Range("N122").Value = -39400

$ws.Range("H126").Value = 2897.4333
$ws.Range("I126").Value = 1339.1052
$ws.Range("J126").Value = 5589.091
$ws.Range("K126").Value = 4017.3156
$ws.Range("L126").Value = 16767.273
$ws.Range("M126").Value = -1547.3156
$ws.Range("N126").Value = -21707.273

$ws.Range("H132").Value = 3915.6765
$ws.Range("I132").Value = 3084.0952
$ws.Range("K132").Value = 9252.285600000001
$ws.Range("M132").Value = -6722.285600000001

$ws.Range("H136").Value = 2595.9482
$ws.Range("I136").Value = 1290.5714
$ws.Range("J136").Value = 4582.391
$ws.Range("K136").Value = 3871.7142
$ws.Range("L136").Value = 13747.173
$ws.Range("M136").Value = -1321.7142
$ws.Range("N136").Value = -18847.173

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 29737.5
$ws.Range("J15").Value = 29737.5
$ws.Range("L15").Value = 29737.5
$ws.Range("N15").Value = -30313.5

$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = $null

$ws.Range("H132").Value = 10755067
$ws.Range("I132").Value = 1809.3334
$ws.Range("K132").Value = 5428.0002
$ws.Range("M132").Value = -2898.0002

$ws.Range("H136").Value = 3084.111
$ws.Range("I136").Value = 1337.3
$ws.Range("J136").Value = 5267.625
$ws.Range("K136").Value = 4011.9
$ws.Range("L136").Value = 15802.875
$ws.Range("M136").Value = -1461.9
$ws.Range("N136").Value = -20902.875
